# US011Steps payment data come from ExcelFile completed
$wb = $excel.ActiveWorkbook

# Add the new worksheet (after the existing last sheet) that will hold
# the campus bank / payment data.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "CampusBankData"

# Populate the payment data.
$newSheet.Range("E10").Value = "4242 4242 4242 4242"
$newSheet.Range("F11").Value = 1223
$newSheet.Range("G12").Value = 123

# Column width for the label column (target stored width 42.5703125 chars;
# the host quantizes ColumnWidth writes to 1/6-character steps, so this is
# the closest input that round-trips to that stored value).
$newSheet.Columns.Item(5).ColumnWidth = 41.666666

# Selection / active state on the new sheet.
$newSheet.Range("G13").Select()
$newSheet.Activate()
